$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the current row 264 (shifting the old
# rows 264-269 down to 267-272), mirroring the rest of the sheet's
# formatting for the inserted rows (style carries from the row below, as
# Excel does on a regular row insert).
$ws.Range("A264:R266").Insert()

# --- New row 264: Alcachofa Española Extra, week of 2021-09-09, Limarí ---
$ws.Range("A264").Value = 9
$ws.Range("B264").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C264").Value = "Metropolitana"
$ws.Range("D264").Value = 44448
$ws.Range("E264").Value = 13
$ws.Range("F264").Value = 100112013
$ws.Range("G264").Value = "Alcachofa"
$ws.Range("H264").Value = "Española"
$ws.Range("I264").Value = "Extra"
$ws.Range("J264").Value = 16
$ws.Range("K264").Value = 15000
$ws.Range("L264").Value = 15000
$ws.Range("M264").Value = 15000
$ws.Range("N264").Value = "`$/caja 25 unidades"
$ws.Range("O264").Value = "Provincia de Limarí"
$ws.Range("P264").Value = 15000
$ws.Range("Q264").Value = 1
$ws.Range("R264").Value = "Hortaliza"

# --- New row 265: Alcachofa Española Primera, week of 2021-09-09, Limarí ---
$ws.Range("A265").Value = 9
$ws.Range("B265").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C265").Value = "Metropolitana"
$ws.Range("D265").Value = 44448
$ws.Range("E265").Value = 13
$ws.Range("F265").Value = 100112013
$ws.Range("G265").Value = "Alcachofa"
$ws.Range("H265").Value = "Española"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 43
$ws.Range("K265").Value = 13000
$ws.Range("L265").Value = 14000
$ws.Range("M265").Value = 13512
$ws.Range("N265").Value = "`$/caja 30 unidades"
$ws.Range("O265").Value = "Provincia de Limarí"
$ws.Range("P265").Value = 450
$ws.Range("Q265").Value = 30
$ws.Range("R265").Value = "Hortaliza"

# --- New row 266: Alcachofa Española Segunda, week of 2021-09-09, Limarí ---
$ws.Range("A266").Value = 9
$ws.Range("B266").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C266").Value = "Metropolitana"
$ws.Range("D266").Value = 44448
$ws.Range("E266").Value = 13
$ws.Range("F266").Value = 100112013
$ws.Range("G266").Value = "Alcachofa"
$ws.Range("H266").Value = "Española"
$ws.Range("I266").Value = "Segunda"
$ws.Range("J266").Value = 25
$ws.Range("K266").Value = 11000
$ws.Range("L266").Value = 12000
$ws.Range("M266").Value = 11480
$ws.Range("N266").Value = "`$/caja 40 unidades"
$ws.Range("O266").Value = "Provincia de Limarí"
$ws.Range("P266").Value = 287
$ws.Range("Q266").Value = 40
$ws.Range("R266").Value = "Hortaliza"
